$d = $word.ActiveDocument
$brk = [char]11

# The original paragraphs each hold one long run of concatenated sentences
# with no separators. The edit breaks each into one <w:t> per sentence/bullet
# joined by manual line breaks (<w:br/>), without disturbing formatting.
# Using Find.Execute to *locate* the target range (no Replace string, so no
# autocorrect/smart-quote mangling is triggered), then assigning the
# multi-line text directly onto the matched Range.Text preserves straight
# apostrophes exactly as in the source.

# --- Edit 1: "Programa resumido" summary paragraph ---
$rng = $d.Content
$found1 = $rng.Find.Execute('1. Técnicas de redação científica, uso de ferramentas de busca, referências bibliográficas e estruturas formais de divulgação científica.2. Desenvolvimento de relatórios científicos.3. Técnicas de coleta, análise e interpretação de dados.4. Importância da revisão por pares e da ética na pesquisa científica.5. Apresentação de resultados de pesquisa de forma clara e eficaz.', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "Edit 1: target paragraph not found" }
$rng.Text = '1. Técnicas de redação científica, uso de ferramentas de busca, referências bibliográficas e estruturas formais de divulgação científica.' + $brk + '2. Desenvolvimento de relatórios científicos.' + $brk + '3. Técnicas de coleta, análise e interpretação de dados.' + $brk + '4. Importância da revisão por pares e da ética na pesquisa científica.' + $brk + '5. Apresentação de resultados de pesquisa de forma clara e eficaz.'
Write-Output "Edit 1 applied: $found1"

# --- Edit 2: "Programa" detailed paragraph ---
$rng = $d.Content
$found2 = $rng.Find.Execute('1. Técnicas de redação científica, uso de ferramentas de busca, referências bibliográficas e estruturas formais de divulgação científica.• Definição e importância do método científico.• Histórico e evolução do método científico.• Aplicação do método científico na Engenharia Química.2. Desenvolvimento de relatórios científicos.• Elaboração de relatórios científicos junto à disciplina Química Geral Experimental.• Redação técnica e clareza na comunicação. 3. Técnicas de coleta, análise e interpretação de dados.• Métodos de coleta de dados.• Ferramentas de software e estatísticas para análise de dados.• Estruturas de artigos científicos. 4. Importância da revisão por pares e da ética na pesquisa científica.• Princípios éticos na pesquisa.• Plágio e integridade científica.• A importância da revisão por pares. 5. Apresentação de resultados de pesquisa de forma clara e eficaz.• Interpretação de resultados experimentais.• Discussão e implicações dos resultados.• Apresentação de resultados em conferências e seminários.• Confecção de pôster e apresentação oral.', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Edit 2: target paragraph not found" }
$rng.Text = '1. Técnicas de redação científica, uso de ferramentas de busca, referências bibliográficas e estruturas formais de divulgação científica.' + $brk + '• Definição e importância do método científico.' + $brk + '• Histórico e evolução do método científico.' + $brk + '• Aplicação do método científico na Engenharia Química.' + $brk + '2. Desenvolvimento de relatórios científicos.' + $brk + '• Elaboração de relatórios científicos junto à disciplina Química Geral Experimental.' + $brk + '• Redação técnica e clareza na comunicação. 3. Técnicas de coleta, análise e interpretação de dados.' + $brk + '• Métodos de coleta de dados.' + $brk + '• Ferramentas de software e estatísticas para análise de dados.' + $brk + '• Estruturas de artigos científicos. 4. Importância da revisão por pares e da ética na pesquisa científica.' + $brk + '• Princípios éticos na pesquisa.' + $brk + '• Plágio e integridade científica.' + $brk + '• A importância da revisão por pares. 5. Apresentação de resultados de pesquisa de forma clara e eficaz.' + $brk + '• Interpretação de resultados experimentais.' + $brk + '• Discussão e implicações dos resultados.' + $brk + '• Apresentação de resultados em conferências e seminários.' + $brk + '• Confecção de pôster e apresentação oral.'
Write-Output "Edit 2 applied: $found2"

# --- Edit 3: "Bibliografia" paragraph ---
$rng = $d.Content
$found3 = $rng.Find.Execute('• LAKATOS, E. M.; MARCONI, M. A. Fundamentos de Metodologia Científica. 7ª ed. São Paulo: Atlas, 2017.• CRESWELL, J. W.; CRESWELL, J. D. Research Design: Qualitative, Quantitative, and Mixed Methods Approaches. 5th ed. Thousand Oaks: Sage Publications, 2018.• MONTGOMERY, D. C. Design and Analysis of Experiments. 9th ed. New York: Wiley, 2019.• RUSSELL, S. W.; MORRISON, D. C. The Grant Application Writer''s Workbook. 2nd ed. Los Angeles: Grant Writers'' Seminars & Workshops LLC, 2018.• ZAR, J. H. Biostatistical Analysis. 5th ed. Upper Saddle River: Pearson Prentice-Hall, 2010.• ALTMAN, D. G. Practical Statistics for Medical Research. 1st ed. Boca Raton: Chapman & Hall/CRC, 1990.', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "Edit 3: target paragraph not found" }
$rng.Text = '• LAKATOS, E. M.; MARCONI, M. A. Fundamentos de Metodologia Científica. 7ª ed. São Paulo: Atlas, 2017.' + $brk + '• CRESWELL, J. W.; CRESWELL, J. D. Research Design: Qualitative, Quantitative, and Mixed Methods Approaches. 5th ed. Thousand Oaks: Sage Publications, 2018.' + $brk + '• MONTGOMERY, D. C. Design and Analysis of Experiments. 9th ed. New York: Wiley, 2019.' + $brk + '• RUSSELL, S. W.; MORRISON, D. C. The Grant Application Writer''s Workbook. 2nd ed. Los Angeles: Grant Writers'' Seminars & Workshops LLC, 2018.' + $brk + '• ZAR, J. H. Biostatistical Analysis. 5th ed. Upper Saddle River: Pearson Prentice-Hall, 2010.' + $brk + '• ALTMAN, D. G. Practical Statistics for Medical Research. 1st ed. Boca Raton: Chapman & Hall/CRC, 1990.'
Write-Output "Edit 3 applied: $found3"
